$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 17:05"

$ws.Range("B4").Value = 4327673
$ws.Range("C4").Value = 11964
$ws.Range("D4").Value = 2061879
$ws.Range("E4").Value = 2116304
$ws.Range("G4").Value = 92
$ws.Range("H4").Value = 149490

$ws.Range("B6").Value = 1424202
$ws.Range("C6").Value = 38708
$ws.Range("D6").Value = 910298
$ws.Range("E6").Value = 481248
$ws.Range("G6").Value = 560
$ws.Range("H6").Value = 32656

$ws.Range("B21").Value = 206363
$ws.Range("C21").Value = 31
$ws.Range("E21").Value = 6561

$ws.Range("B40").Value = 62908
$ws.Range("C40").Value = 2012
$ws.Range("D40").Value = 28603
$ws.Range("E40").Value = 33242
$ws.Range("G40").Value = 8
$ws.Range("H40").Value = 1063

$ws.Range("B41").Value = 61388
$ws.Range("C41").Value = 710
$ws.Range("D41").Value = 26959
$ws.Range("E41").Value = 33965
$ws.Range("G41").Value = 7
$ws.Range("H41").Value = 464

$ws.Range("D45").Value = 45521
$ws.Range("E45").Value = 4821

$ws.Range("B63").Value = 23034
$ws.Range("C63").Value = 206
$ws.Range("E63").Value = 6390
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 735

$ws.Range("B114").Value = 2777
$ws.Range("C114").Value = 7
$ws.Range("E114").Value = 660

$ws.Range("B118").Value = 2510
$ws.Range("C118").Value = 7
$ws.Range("D118").Value = 1911
$ws.Range("E118").Value = 476

$ws.Range("B141").Value = 1168
$ws.Range("C141").Value = 14
$ws.Range("D141").Value = 1041
$ws.Range("E141").Value = 116
$ws.Range("H141").Value = 11

$ws.Range("B142").Value = 1162
$ws.Range("C142").Value = 7
$ws.Range("D142").Value = 641
$ws.Range("E142").Value = 449
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 72
